$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 23

$ws.Cells.Item($row, 1).Value = "YJZTXE"
$ws.Cells.Item($row, 2).Value = "Cinta Flex Epson"
$ws.Cells.Item($row, 3).Value = "LX 300+"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 70000
$ws.Cells.Item($row, 6).Value = 5
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E23-D23)*G23"
$ws.Cells.Item($row, 9).Formula = "=D23*F23"
$ws.Cells.Item($row, 10).Value = 0
